# phonegap project.xlsx - apply commit:
#   > add array sorting to unify order
#   > adjust alert msg format
#   > add custom confirm dialog
#   > update default open list rule

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "temp note": two new rows documenting the jQuery Mobile confirm
# dialog research (A23 bold heading, A24 the stackoverflow reference).
# ---------------------------------------------------------------------
$tn = $wb.Worksheets.Item("temp note")

# Match A24's formatting/order of creation first so shared-string order
# comes out the same as the source workbook (url string before caption).
$tn.Range("A14").Copy()
$tn.Range("A24").PasteSpecial(-4122)   # xlPasteFormats
$tn.Range("A24").Value = "http://stackoverflow.com/questions/5747382/jquery-mobile-alert-confirmation-dialogs"

$tn.Range("A12").Copy()
$tn.Range("A23").PasteSpecial(-4122)   # xlPasteFormats
$tn.Range("A23").Value = "jQuery mobile confirm dialog"

$tn.Activate() | Out-Null
$tn.Range("A23").Select() | Out-Null

# ---------------------------------------------------------------------
# Sheet "decision": unify the "Button function redesign" / "Draw" blocks
# (remove the stray blank separator rows above them so everything sits
# together), add a red custom-confirm-dialog caveat note, and drop the
# obsolete "open list : SESSION change to ??" / "constant by javascript"
# scratch notes entirely.
# ---------------------------------------------------------------------
$dec = $wb.Worksheets.Item("decision")

# Rows 15-18 were just a lone empty A15 plus blank rows -- delete them so
# the "Button function redesign" block (old rows 19-25) shifts up to 15-21.
$dec.Rows("15:18").Delete()

# New note (row 22, which is now blank) in red text.
$dec.Range("B22").Value = "custom confirm msg can't reload page after confirm"
$dec.Range("B22").Font.Name = "Calibri"
$dec.Range("B22").Font.Size = 12
$dec.Range("B22").Font.Color = 255

# Drop the old "open list" (was row 34) and "constant by javascript" (was
# row 39) scratch blocks -- now sitting at rows 26-37 after the shift above.
$dec.Rows("26:37").Delete()

$dec.Activate() | Out-Null
$dec.Range("A30").Select() | Out-Null

Write-Host "phonegap project.xlsx updated"
